$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $newValue) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $newValue
    $c.NumberFormat = "General"
    $c.Style = "Normal"
}

Set-TextValue "D2" "66.835.89"
Set-TextValue "E2" "  -3.68%  "

Set-TextValue "D3" "3.526.78"
Set-TextValue "E3" "  -4.02%  "

Set-TextValue "E4" "  +0.00%  "

Set-TextValue "D5" "605.94"
Set-TextValue "E5" "  -5.49%  "

Set-TextValue "D6" "154.20"
Set-TextValue "E6" "  -3.43%  "

Set-TextValue "D7" "3.525.56"
Set-TextValue "E7" "  -3.92%  "

Set-TextValue "E8" "  +0.12%  "

Set-TextValue "D9" "0.485"
Set-TextValue "E9" "  -2.51%  "

Set-TextValue "E10" "  -2.22%  "

Set-TextValue "D11" "6.83"
Set-TextValue "E11" "  -3.51%  "

Set-TextValue "D12" "0.431"
Set-TextValue "E12" "  -3.89%  "

Set-TextValue "D13" "0.0000221"
Set-TextValue "E13" "  -4.45%  "

Set-TextValue "D14" "4.128.79"
Set-TextValue "E14" "  -3.87%  "

Set-TextValue "D15" "31.91"
Set-TextValue "E15" "  -2.29%  "

Set-TextValue "D16" "3.520.25"
Set-TextValue "E16" "  -4.33%  "

Set-TextValue "D17" "66.893.53"
Set-TextValue "E17" "  -3.62%  "

Set-TextValue "E18" "  +0.83%  "

Set-TextValue "D19" "6.36"
Set-TextValue "E19" "  -1.85%  "

Set-TextValue "D20" "15.41"
Set-TextValue "E20" "  -3.49%  "

Set-TextValue "D21" "451.10"
Set-TextValue "E21" "  -3.08%  "

Set-TextValue "D22" "9.37"
Set-TextValue "E22" "  -5.29%  "

Set-TextValue "D23" "0.638"
Set-TextValue "E23" "  -1.34%  "

Set-TextValue "D24" "79.11"
Set-TextValue "E24" "  -0.33%  "

Set-TextValue "D25" "3.672.87"
Set-TextValue "E25" "  -3.89%  "

Set-TextValue "E26" "  +0.05%  "

Set-TextValue "D27" "0.0000123"
Set-TextValue "E27" "  -2.02%  "

Set-TextValue "D28" "10.25"
Set-TextValue "E28" "  -5.74%  "

Set-TextValue "D29" "8.30"
Set-TextValue "E29" "  -8.13%  "

Set-TextValue "D30" "2.55"
Set-TextValue "E30" "  -2.74%  "

Set-TextValue "E31" "  -1.60%  "

Set-TextValue "E32" "  +0.13%  "

Set-TextValue "D33" "25.90"
Set-TextValue "E33" "  -3.49%  "

Set-TextValue "D34" "1.89"
Set-TextValue "E34" "  -5.51%  "

Set-TextValue "D35" "6.18"
Set-TextValue "E35" "  -4.22%  "

Set-TextValue "D36" "0.157"
Set-TextValue "E36" "  -4.83%  "

Set-TextValue "D37" "3.524.42"
Set-TextValue "E37" "  -3.94%  "

Set-TextValue "D38" "8.08"
Set-TextValue "E38" "  -4.14%  "

Set-TextValue "E39" "  +0.00%  "

Set-TextValue "E40" "  +0.03%  "

Set-TextValue "D41" "176.77"
Set-TextValue "E41" "  -0.87%  "

Set-TextValue "D42" "5.60"
Set-TextValue "E42" "  -4.86%  "

Set-TextValue "B43" "Hedera"
Set-TextValue "C43" "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue "D43" "0.0876"
Set-TextValue "E43" "  -2.63%  "

Set-TextValue "B44" "Stacks"
Set-TextValue "C44" "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-TextValue "D44" "2.13"
Set-TextValue "E44" "  -2.48%  "

Set-TextValue "D45" "0.891"
Set-TextValue "E45" "  -3.66%  "

Set-TextValue "D46" "45.71"
Set-TextValue "E46" "  -2.04%  "

Set-TextValue "D47" "28.46"
Set-TextValue "E47" "  +4.47%  "

Set-TextValue "D48" "2.69"
Set-TextValue "E48" "  -1.70%  "

Set-TextValue "D49" "1.22"
Set-TextValue "E49" "  -1.97%  "

Set-TextValue "D50" "7.64"
Set-TextValue "E50" "  -2.42%  "

Set-TextValue "E51" "  -3.23%  "
